$d = $word.ActiveDocument

# 1. Expand the title "MAKE" -> "MAKE Technical Documentation".
#    Scope the Find/Replace to the Title paragraph only so that the many
#    other "MAKE" occurrences throughout the body are left untouched.
$titlePara = $d.Paragraphs(1)
$titlePara.Range.Find.Execute("MAKE", $true, $true, $false, $false, $false, `
                               $true, 1, $false, "MAKE Technical Documentation", 2)

# 2. Fix the typo in the Abstract paragraph: "resident" -> "recognized".
#    Scope the Find/Replace to the Abstract paragraph only (and match the
#    whole word) so the single occurrence is corrected without touching
#    any other text.
$abstractPara = $d.Paragraphs(4)
$abstractPara.Range.Find.Execute("resident", $true, $true, $false, $false, $false, `
                                  $true, 1, $false, "recognized", 2)
